$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '54.133.08'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -8.47%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.423.93'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -14.44%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '463.12'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -7.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.33'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.75%  '
$ws.Range('E7').Value = '  -0.62%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.490'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -7.31%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.432.01'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -14.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0947'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -8.47%  '
$ws.Range('E11').Value = '  -9.33%  '
$ws.Range('E12').Value = '  -8.14%  '
$ws.Range('E13').Value = '  -4.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.832.96'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -15.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '53.993.26'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -8.91%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.68'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -9.64%  '
$ws.Range('E17').Value = '  -2.82%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.432.96'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -14.03%  '
$ws.Range('E19').Value = '  -10.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '312.79'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -11.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.38'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -14.97%  '
$ws.Range('E22').Value = '  +0.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.67'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.59%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.39'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -13.57%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '56.51'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -10.37%  '
$ws.Range('E26').Value = '  +1.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.384'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -10.47%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.507.08'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -15.38%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.153'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -10.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.18'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.996'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.34%  '
$ws.Range('E32').Value = '  -10.94%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '149.61'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.67'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -7.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.40'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -12.55%  '
$ws.Range('E36').Value = '  -5.66%  '
$ws.Range('E37').Value = '  -15.20%  '
$ws.Range('E38').Value = '  -7.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.799'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -11.53%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '33.51'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -8.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.992'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.70%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.605'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.76%  '
$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0529'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.26%  '
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.29'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.22%  '
$ws.Range('B45').Value = 'WhiteBITCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.15'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.90%  '
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.24'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -8.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.971.24'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -11.02%  '
$ws.Range('E48').Value = '  -2.23%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0866'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.36'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.64%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.53'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -14.75%  '
